$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume-1h (E) columns for rows 2-49
# Row 2
$ws.Range("D2").Value = "'28.247.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +3.06%  "
$ws.Range("E2").ClearFormats()
# Row 3
$ws.Range("D3").Value = "'1.811.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +4.16%  "
$ws.Range("E3").ClearFormats()
# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.19%  "
$ws.Range("E4").ClearFormats()
# Row 5
$ws.Range("D5").Value = "'326.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.62%  "
$ws.Range("E5").ClearFormats()
# Row 6
$ws.Range("D6").Value = "'0.9997"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  -0.06%  "
$ws.Range("E6").ClearFormats()
# Row 7
$ws.Range("D7").Value = "'0.4355"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +2.91%  "
$ws.Range("E7").ClearFormats()
# Row 8
$ws.Range("D8").Value = "'0.3660"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +2.09%  "
$ws.Range("E8").ClearFormats()
# Row 9
$ws.Range("D9").Value = "'44.92"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -1.01%  "
$ws.Range("E9").ClearFormats()
# Row 10
$ws.Range("D10").Value = "'0.07670"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +3.66%  "
$ws.Range("E10").ClearFormats()
# Row 11
$ws.Range("D11").Value = "'1.142"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  +2.93%  "
$ws.Range("E11").ClearFormats()
# Row 12
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  -0.15%  "
$ws.Range("E12").ClearFormats()
# Row 13
$ws.Range("D13").Value = "'22.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +3.29%  "
$ws.Range("E13").ClearFormats()
# Row 14
$ws.Range("D14").Value = "'6.295"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +3.66%  "
$ws.Range("E14").ClearFormats()
# Row 15
$ws.Range("D15").Value = "'7.523"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  +4.87%  "
$ws.Range("E15").ClearFormats()
# Row 16
$ws.Range("D16").Value = "'1.823.77"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +4.83%  "
$ws.Range("E16").ClearFormats()
# Row 17
$ws.Range("D17").Value = "'95.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +9.53%  "
$ws.Range("E17").ClearFormats()
# Row 18
$ws.Range("D18").Value = "'0.00001082"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.76%  "
$ws.Range("E18").ClearFormats()
# Row 19
$ws.Range("D19").Value = "'0.06534"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +4.57%  "
$ws.Range("E19").ClearFormats()
# Row 20
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.02%  "
$ws.Range("E20").ClearFormats()
# Row 21
$ws.Range("D21").Value = "'17.46"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +3.88%  "
$ws.Range("E21").ClearFormats()
# Row 22
$ws.Range("D22").Value = "'6.238"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +2.43%  "
$ws.Range("E22").ClearFormats()
# Row 23
$ws.Range("D23").Value = "'28.262.27"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +2.99%  "
$ws.Range("E23").ClearFormats()
# Row 24
$ws.Range("D24").Value = "'11.56"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.20%  "
$ws.Range("E24").ClearFormats()
# Row 25
$ws.Range("D25").Value = "'2.080"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -10.23%  "
$ws.Range("E25").ClearFormats()
# Row 26
$ws.Range("D26").Value = "'161.44"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +6.90%  "
$ws.Range("E26").ClearFormats()
# Row 27
$ws.Range("D27").Value = "'20.67"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +1.58%  "
$ws.Range("E27").ClearFormats()
# Row 28
$ws.Range("D28").Value = "'2.026.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +4.61%  "
$ws.Range("E28").ClearFormats()
# Row 29
$ws.Range("D29").Value = "'2.287"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -2.21%  "
$ws.Range("E29").ClearFormats()
# Row 30
$ws.Range("D30").Value = "'128.78"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.95%  "
$ws.Range("E30").ClearFormats()
# Row 31
$ws.Range("D31").Value = "'1.209"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  -0.31%  "
$ws.Range("E31").ClearFormats()
# Row 32
$ws.Range("D32").Value = "'5.937"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  +4.91%  "
$ws.Range("E32").ClearFormats()
# Row 33
$ws.Range("D33").Value = "'0.09174"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +0.46%  "
$ws.Range("E33").ClearFormats()
# Row 34
$ws.Range("D34").Value = "'3.460"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -5.69%  "
$ws.Range("E34").ClearFormats()
# Row 35
$ws.Range("D35").Value = "'12.99"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +2.82%  "
$ws.Range("E35").ClearFormats()
# Row 36
$ws.Range("D36").Value = "'0.02350"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  +2.90%  "
$ws.Range("E36").ClearFormats()
# Row 37
$ws.Range("D37").Value = "'0.2172"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  +2.13%  "
$ws.Range("E37").ClearFormats()
# Row 38
$ws.Range("D38").Value = "'5.191"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +2.50%  "
$ws.Range("E38").ClearFormats()
# Row 39
$ws.Range("D39").Value = "'0.6570"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  +2.94%  "
$ws.Range("E39").ClearFormats()
# Row 40
$ws.Range("D40").Value = "'0.06205"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  +2.39%  "
$ws.Range("E40").ClearFormats()
# Row 41
$ws.Range("D41").Value = "'1.194"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -0.14%  "
$ws.Range("E41").ClearFormats()
# Row 42
$ws.Range("D42").Value = "'8.127"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +3.13%  "
$ws.Range("E42").ClearFormats()
# Row 43
$ws.Range("D43").Value = "'1.424"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  +0.51%  "
$ws.Range("E43").ClearFormats()
# Row 44
$ws.Range("D44").Value = "'0.9997"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -0.04%  "
$ws.Range("E44").ClearFormats()
# Row 45
$ws.Range("D45").Value = "'13.81"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +1.42%  "
$ws.Range("E45").ClearFormats()
# Row 46
$ws.Range("D46").Value = "'0.6103"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +4.10%  "
$ws.Range("E46").ClearFormats()
# Row 47
$ws.Range("D47").Value = "'3.743"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +0.65%  "
$ws.Range("E47").ClearFormats()
# Row 48
$ws.Range("D48").Value = "'125.82"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +0.70%  "
$ws.Range("E48").ClearFormats()
# Row 49
$ws.Range("D49").Value = "'2.016"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +3.49%  "
$ws.Range("E49").ClearFormats()

# Row 50 and 51: Cronos/EOS swap with updated price/volume
# Row 50 becomes EOS
$ws.Range("B50").Value = "'EOS"
$ws.Range("B50").ClearFormats()
$ws.Range("C50").Value = "'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("C50").ClearFormats()
$ws.Range("D50").Value = "'1.156"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +3.31%  "
$ws.Range("E50").ClearFormats()

# Row 51 becomes Cronos
$ws.Range("B51").Value = "'Cronos"
$ws.Range("B51").ClearFormats()
$ws.Range("C51").Value = "'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("C51").ClearFormats()
$ws.Range("D51").Value = "'0.07001"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +2.28%  "
$ws.Range("E51").ClearFormats()
